$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Overview" (sheet1): rows for the two source files swap order
# (f88f1c75 now handed back -> row 2, 551ef9d2 still pending -> row 3)
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "f88f1c75-83c7-44b2-be9d-341104fcb25d.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-03-24 14:48:23"

$ws1.Range("A3").Value = "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-03-24 14:47:32"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a387dd077c39973a1e1ecc886b6ad9d2393b080b/e2e/f88f1c75-83c7-44b2-be9d-341104fcb25d.md", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3407373e5ef3d9cd09a1eba61467cdb0041a0c9e/e2e/551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md", "", "", "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md")

# -----------------------------------------------------------------
# Sheet "zh-cn" (sheet2): f88f1c75 row handed back (gets target/handback
# file + datetime, new Status), rows reordered (f88f1c75 -> row2,
# 551ef9d2 -> row3)
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "f88f1c75-83c7-44b2-be9d-341104fcb25d.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-24 14:48:18"
$ws2.Range("F2").Value = "f88f1c75-83c7-44b2-be9d-341104fcb25d.md"
$ws2.Range("G2").Value = "f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-24 14:48:48"
$ws2.Range("J2").Value = "Include"

$ws2.Range("A3").Value = "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.460871851b3c3f69f4cdb5f568904f746a515d1c.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-24 14:47:27"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("J3").Value = "Include"

$ws2.Range("F2").Style = "HyperLink"
$ws2.Range("G2").Style = "HyperLink"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a387dd077c39973a1e1ecc886b6ad9d2393b080b/e2e/f88f1c75-83c7-44b2-be9d-341104fcb25d.md", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a3384ba76bdf6b091b306fb3c250dab4e483a9a5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.zh-cn.xlf", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/a387dd077c39973a1e1ecc886b6ad9d2393b080b/e2e/f88f1c75-83c7-44b2-be9d-341104fcb25d.md", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a3384ba76bdf6b091b306fb3c250dab4e483a9a5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.zh-cn.xlf", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3407373e5ef3d9cd09a1eba61467cdb0041a0c9e/e2e/551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md", "", "", "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e148cd17d31df31b4a8823697e7085eaa0ed62ef/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/551ef9d2-4f53-4ab1-9a3a-b25095c949aa.460871851b3c3f69f4cdb5f568904f746a515d1c.zh-cn.xlf", "", "", "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.460871851b3c3f69f4cdb5f568904f746a515d1c.zh-cn.xlf")

# -----------------------------------------------------------------
# Sheet "de-de" (sheet3): same pattern as zh-cn, using de-de xlf files
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "f88f1c75-83c7-44b2-be9d-341104fcb25d.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-24 14:48:23"
$ws3.Range("F2").Value = "f88f1c75-83c7-44b2-be9d-341104fcb25d.md"
$ws3.Range("G2").Value = "f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-24 14:48:54"
$ws3.Range("J2").Value = "Include"

$ws3.Range("A3").Value = "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.460871851b3c3f69f4cdb5f568904f746a515d1c.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-24 14:47:32"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("J3").Value = "Include"

$ws3.Range("F2").Style = "HyperLink"
$ws3.Range("G2").Style = "HyperLink"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a387dd077c39973a1e1ecc886b6ad9d2393b080b/e2e/f88f1c75-83c7-44b2-be9d-341104fcb25d.md", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f5e5c43f3582d1630d5eb7ab1143cd42aeca3b62/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.de-de.xlf", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/a387dd077c39973a1e1ecc886b6ad9d2393b080b/e2e/f88f1c75-83c7-44b2-be9d-341104fcb25d.md", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f5e5c43f3582d1630d5eb7ab1143cd42aeca3b62/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.de-de.xlf", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3407373e5ef3d9cd09a1eba61467cdb0041a0c9e/e2e/551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md", "", "", "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1ddce9b2df8afb83fc283938c526c02aba70f64c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/551ef9d2-4f53-4ab1-9a3a-b25095c949aa.460871851b3c3f69f4cdb5f568904f746a515d1c.de-de.xlf", "", "", "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.460871851b3c3f69f4cdb5f568904f746a515d1c.de-de.xlf")

$wb.Save()
